# DBinfo.xlsx edit: update "DB Updated Date" value and fix the saved selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 holds the text "2021.02.02" (DB Updated Date). Update it to "2021.02.18".
# Assign via a literal-text formula and paste back as a value so Excel doesn't
# reinterpret the dotted string as a date (and so no extra cell style gets
# introduced for the cell).
$ws.Range("B1").Formula = "=""2021.02.18"""
$ws.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# The saved sheet view should have its selection/active cell on B1 (it was H1).
$ws.Range("B1").Select()
